# Updated cryptos list (price + volume(1h) refresh, plus Aave/FraxShare row swap)
# Commit: Updated cryptos list on Fri Jan 12 11:47:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format before assigning, so numeric-looking strings (e.g. "310.48")
# are preserved verbatim as text, matching the original inlineStr cell type,
# instead of being auto-converted into floating point numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.093.10'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.649.90'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.48'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.78'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.01%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.76'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.32'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.11'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.051.05'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.655.18'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '46.135.90'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.11%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.88'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '282.19'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.07'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '30.20'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.60%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.60'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '38.69'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.27'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.77'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.36'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.38'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0843'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.33%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '15.84'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.07%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.66%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.155.44'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.21%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '93.98'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.27'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '111.20'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.902.59'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.13%  '
